$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62; existing rows 62-118 shift down to 63-119.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new record.
$ws.Cells.Item(62, 1).Value2 = 8
$ws.Cells.Item(62, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(62, 3).Value2 = "Coquimbo"
$ws.Cells.Item(62, 4).Value2 = 44673
$ws.Cells.Item(62, 5).Value2 = 4
$ws.Cells.Item(62, 6).Value2 = "Fruta"
$ws.Cells.Item(62, 7).Value2 = 100109
$ws.Cells.Item(62, 8).Value2 = "Uva"
$ws.Cells.Item(62, 9).Value2 = 100109001
$ws.Cells.Item(62, 10).Value2 = "Uva"
$ws.Cells.Item(62, 11).Value2 = "Red Globe"
$ws.Cells.Item(62, 12).Value2 = "Primera"
$ws.Cells.Item(62, 13).Value2 = 400
$ws.Cells.Item(62, 14).Value2 = 9000
$ws.Cells.Item(62, 15).Value2 = 9500
$ws.Cells.Item(62, 16).Value2 = 9250
$ws.Cells.Item(62, 17).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(62, 18).Value2 = "Provincia del Elquí"
$ws.Cells.Item(62, 19).Value2 = 514
$ws.Cells.Item(62, 20).Value2 = 18
